$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58, shifting existing rows 58:69 down to 59:70
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the new weekly record
$ws.Cells.Item(58, 1).Value = 5
$ws.Cells.Item(58, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(58, 3).Value = "Maule"
$ws.Cells.Item(58, 4).Value = 44551
$ws.Cells.Item(58, 5).Value = 7
$ws.Cells.Item(58, 6).Value = 100112001
$ws.Cells.Item(58, 7).Value = "Berenjena"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 100
$ws.Cells.Item(58, 11).Value = 10000
$ws.Cells.Item(58, 12).Value = 10000
$ws.Cells.Item(58, 13).Value = 10000
$ws.Cells.Item(58, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(58, 15).Value = "Región del Maule"
$ws.Cells.Item(58, 16).Value = 200
$ws.Cells.Item(58, 17).Value = 50
$ws.Cells.Item(58, 18).Value = "Hortaliza"
